$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number need to be forced to Text
# format first, otherwise Excel auto-converts the assigned string into a
# numeric value (the source file stores every Price/Volume cell as text).
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '27.625.03'
$ws.Range('E2').Value = '  -1.03%  '
$ws.Range('D3').Value = '1.878.49'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '331.49'
$ws.Range('E5').Value = '  +1.44%  '
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('E7').Value = '  +3.03%  '
$ws.Range('D8').Value = '0.3963'
$ws.Range('D9').Value = '47.80'
$ws.Range('E9').Value = '  -3.12%  '
$ws.Range('D10').Value = '0.08010'
$ws.Range('E10').Value = '  -2.45%  '
$ws.Range('D11').Value = '1.023'
$ws.Range('E11').Value = '  -1.15%  '
$ws.Range('D12').Value = '21.73'
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('D13').Value = '1.877.30'
$ws.Range('E13').Value = '  -1.80%  '
$ws.Range('D14').Value = '5.957'
$ws.Range('E14').Value = '  -0.05%  '
$ws.Range('D15').Value = '7.154'
$ws.Range('E15').Value = '  -1.91%  '
$ws.Range('D16').Value = '1.004'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').Value = '86.89'
$ws.Range('E17').Value = '  -2.08%  '
$ws.Range('D18').Value = '0.00001041'
$ws.Range('E18').Value = '  -1.21%  '
$ws.Range('D19').Value = '0.06618'
$ws.Range('E19').Value = '  +0.71%  '
$ws.Range('D20').Value = '17.22'
$ws.Range('E20').Value = '  -1.01%  '
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('D22').Value = '27.648.69'
$ws.Range('E22').Value = '  -1.01%  '
$ws.Range('D23').Value = '5.493'
$ws.Range('E23').Value = '  -2.39%  '
$ws.Range('E24').Value = '  -0.60%  '
$ws.Range('D25').Value = '2.301'
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('D26').Value = '2.099.23'
$ws.Range('E26').Value = '  -1.66%  '
$ws.Range('D27').Value = '156.15'
$ws.Range('E27').Value = '  +1.34%  '
$ws.Range('D28').Value = '20.22'
$ws.Range('E28').Value = '  +1.89%  '
$ws.Range('D29').Value = '2.086'
$ws.Range('E29').Value = '  -0.64%  '
$ws.Range('D30').Value = '5.550'
$ws.Range('E30').Value = '  -2.43%  '
$ws.Range('D31').Value = '122.11'
$ws.Range('E31').Value = '  -0.81%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = '0.09540'
$ws.Range('E32').Value = '  +0.18%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = '0.9623'
$ws.Range('E33').Value = '  +0.91%  '
$ws.Range('D34').Value = '1.453'
$ws.Range('E34').Value = '  -1.14%  '
$ws.Range('D35').Value = '3.631'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('D36').Value = '5.293'
$ws.Range('E36').Value = '  -2.76%  '
$ws.Range('D37').Value = '0.06107'
$ws.Range('E37').Value = '  +0.29%  '
$ws.Range('D38').Value = '0.02250'
$ws.Range('E38').Value = '  -1.25%  '
$ws.Range('D39').Value = '1.225'
$ws.Range('E39').Value = '  -1.76%  '
$ws.Range('D40').Value = '8.099'
$ws.Range('E40').Value = '  -5.24%  '
$ws.Range('D42').Value = '0.5986'
$ws.Range('E42').Value = '  -1.75%  '
$ws.Range('D43').Value = '0.1895'
$ws.Range('E43').Value = '  +0.17%  '
$ws.Range('D44').Value = '10.24'
$ws.Range('E44').Value = '  -4.21%  '
$ws.Range('D45').Value = '1.251'
$ws.Range('E45').Value = '  -3.96%  '
$ws.Range('D46').Value = '0.5684'
$ws.Range('E46').Value = '  -1.99%  '
$ws.Range('D47').Value = '12.22'
$ws.Range('E47').Value = '  -3.18%  '
$ws.Range('D48').Value = '3.408'
$ws.Range('E48').Value = '  -0.41%  '
$ws.Range('D49').Value = '1.930'
$ws.Range('E49').Value = '  -2.69%  '
$ws.Range('E50').Value = '  -1.13%  '
$ws.Range('D51').Value = '110.96'
$ws.Range('E51').Value = '  +0.77%  '
